$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = "Git werkt en zelf het digitale github aan de gang gekregen dus ga voortaan vanuit daar werken."
$ws.Range("B12").Value = "Mocht ik ooit weer offline gaan werken kan ik gewoon alle bestanden downloaden en in de lokale map vervangen."
$ws.Range("B13").Value = "Nu is het plan om in het tabblad: `"BSCity Uitgewerkt`" de inhoud van het spel uit te werken."
$ws.Range("B14").Value = "Dit mag al alles bevatten wat er ooit in het spel gaat komen maar het moet ook duidelijk zijn wat er al in zit, wat er al volgende toegevoegd moet worden en wat functies zijn die echt voor over een lange tijd zijn."
$ws.Range("B16").Value = "Verder heb ik geëxperimenteerd met dit excel bestand toevoegen aan github. Ik kan het daarin niet bewerken maar ga dit nieuwe stukje tekst vanaf regel 11 nu uploaden om te kijken of ik het kan overschrijven en niet dat er dadelijk allemaal bestanden met dezelfde naam in zitten."

$ws.Range("B17").Select()
